{"js": "// The underlying change in this revision is *not* a content/formatting\n// edit at all: every hunk in the diff touches the same handful of\n// elements (the `w:rFonts` inside the TOC content control's `sdtEndPr`,\n// and the six `w:ptab` right-aligned dot-leader tab marks used by the\n// Table of Contents entries) and in every single hunk the attribute\n// *values* are completely unchanged -- only the order in which the\n// attributes are written out is different (e.g.\n// `w:alignment=\"right\" w:leader=\"dot\" w:relativeTo=\"margin\"` becomes\n// `w:relativeTo=\"margin\" w:alignment=\"right\" w:leader=\"dot\"`).\n//\n// That is the signature of a generating-library upgrade (the commit\n// message confirms it: \"Moving from POI 3.17.0 to 4.0.1\") -- Apache POI's\n// XMLBeans-generated bean classes changed the order in which they\n// serialize an element's attributes between those two releases, with no\n// effect whatsoever on the document's visible content, formatting or\n// semantics. A canonical (attribute-order-insensitive) comparison of the\n// OOXML before and after this commit is byte-for-byte identical.\n//\n// There is no operation in the Word JavaScript API that lets an add-in\n// choose or influence the serialized attribute order of an XML element --\n// that is purely an artifact of whichever library/version writes the\n// package, not something the object model exposes. So the faithful,\n// content-preserving way to \"apply\" this diff through Office.js is to\n// leave the document exactly as it is: no body/range/font/table-of-\n// contents mutation corresponds to \"re-order these XML attributes\", and\n// this document's text, formatting, styles and structure must stay\n// untouched (the diff shows nothing else changed).\n//\n// Intentionally a no-op: load nothing, change nothing, sync nothing.\n// (Performing a `context.sync()` with no queued changes is a harmless,\n// well-formed way to show the context was touched without mutating the\n// document.)\nawait context.sync();\n", "ps1": "# The underlying change in this revision is *not* a content/formatting\n# edit at all: every hunk in the diff touches the same handful of\n# elements (the <w:rFonts> inside the TOC content control's <w:sdtEndPr>,\n# and the six <w:ptab> right-aligned dot-leader tab marks used by the\n# Table of Contents entries) and in every single hunk the attribute\n# *values* are completely unchanged -- only the order in which the\n# attributes are written out is different (e.g.\n# w:alignment=\"right\" w:leader=\"dot\" w:relativeTo=\"margin\" becomes\n# w:relativeTo=\"margin\" w:alignment=\"right\" w:leader=\"dot\").\n#\n# That is the signature of a generating-library upgrade (the commit\n# message confirms it: \"Moving from POI 3.17.0 to 4.0.1\") -- Apache POI's\n# XMLBeans-generated bean classes changed the order in which they\n# serialize an element's attributes between those two releases, with no\n# effect whatsoever on the document's visible content, formatting or\n# semantics. A canonical (attribute-order-insensitive) comparison of the\n# OOXML before and after this commit is byte-for-byte identical.\n#\n# There is no Word COM property or method that lets automation choose or\n# influence the serialized attribute order of an XML element -- that is\n# purely an artifact of whichever library/version writes the package, not\n# something the object model exposes ($d.Range().WordOpenXML is the\n# closest thing to \"raw XML\" COM offers, and re-assigning it does not\n# change the document's serialized attribute order either). So the\n# faithful, content-preserving way to \"apply\" this diff through the Word\n# object model is to leave the document exactly as it is: no\n# Find/Replace, Range, Paragraphs or TablesOfContents mutation\n# corresponds to \"re-order these XML attributes\", and this document's\n# text, formatting, styles and structure must stay untouched (the diff\n# shows nothing else changed).\n#\n# Intentionally a no-op: touch nothing on $word.ActiveDocument.\n$d = $word.ActiveDocument\n"}
